$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new value "SmartTube" in column A, row 21 (next empty row)
$ws.Range("A21").Value = "SmartTube"
